$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from SCD0266 to SCD0016
$ws.Name = "SCD0016"

# 2. Update TC_ID cells from DGS-281 to SCD0016-040
$ws.Range("B2").Value = "SCD0016-040"
$ws.Range("B3").Value = "SCD0016-040"
$ws.Range("B4").Value = "SCD0016-040"
$ws.Range("B5").Value = "SCD0016-040"

# 3. Apply horizontal-left alignment to the used data range A1:P5
$ws.Range("A1:P5").HorizontalAlignment = -4131  # xlLeft

# 4. Update sheet view: scroll position and selection
$ws.Range("B6").Select()
